$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the Python-list-formatted rail sequence strings in column D
# with plain concatenated letter strings (no brackets/quotes/commas).

$ws.Range("D2").Value  = "ABCDDCBACBADBDCAADBCCADBDABCCDBADBACABDCBCADDABC"
$ws.Range("D3").Value  = "ABCDDCBACBADBDCAADBCCADBDABCCDBADBACABDCBCADDABC"
$ws.Range("D4").Value  = "DABCABDCCBADDBCABCDABDACCABDDABCCBADCABDCDABBDAC"
$ws.Range("D5").Value  = "DABCABDCCBADDBCABCDABDACCABDDABCCBADCABDCDABBDAC"
$ws.Range("D6").Value  = "CBADACDBABDCBACDABDCBDACBDCADBACABDCCDBABCDABDAC"
$ws.Range("D7").Value  = "CBADACDBABDCBACDABDCBDACBDCADBACABDCCDBABCDABDAC"
$ws.Range("D8").Value  = "BDACABCDBCADACBDCDABDBCACABDADCBABDCBACDBCDADABC"
$ws.Range("D9").Value  = "BDACABCDBCADACBDCDABDBCACABDADCBABDCBACDBCDADABC"
$ws.Range("D10").Value = "DBCACBDAACBDACDBDCABACDBADCBDABCDBACCDBAABCDADCB"
$ws.Range("D11").Value = "DBCACBDAACBDACDBDCABACDBADCBDABCDBACCDBAABCDADCB"
